# Auto-committed on 2023/06/30 週五 17:05:28.55
#
# Changes applied:
#  - Sheet "DBD" (sheet1):
#      * D1: "火險初保檔" -> "擔保品火險檔"
#      * G24: clear the old remark text (content removed, formatting kept)
#      * Selection moves from H22 to I23
#  - Sheet "DBS" (sheet2):
#      * Two new lookup rows appended (findOrigInsuNoAll, findOrigInsuNoFirst)
#      * Column A widened to fit the longer function names
#      * Selection moves to the new last row C7

$wb = $excel.ActiveWorkbook

# ---- Sheet "DBD" ----
$dbd = $wb.Worksheets.Item("DBD")

$dbd.Range("D1").Value = "擔保品火險檔"
$dbd.Range("G24").ClearContents()

$dbd.Activate()
$dbd.Range("I23").Select()

# ---- Sheet "DBS" ----
$dbs = $wb.Worksheets.Item("DBS")

$dbs.Range("A6").Value = "findOrigInsuNoAll"
$dbs.Range("B6").Value = "OrigInsuNo ="
$dbs.Range("C6").Value = "ClCode1,ClCode2,ClNo"

$dbs.Range("A7").Value = "findOrigInsuNoFirst"
$dbs.Range("B7").Value = "ClCode1 = ,AND ClCode2 = ,AND ClNo = ,AND OrigInsuNo ="
$dbs.Range("C7").Value = "CreateDate ASC"

$dbs.Columns.Item(1).ColumnWidth = 163/7

$dbs.Activate()
$dbs.Range("C7").Select()
